$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new value would otherwise be
# auto-parsed as a number by Excel (column D stores prices as text).
$textFormatCells = @("D5", "D8", "D10", "D15", "D16", "D18", "D19", "D25", "D27", "D30", "D36", "D37", "D39", "D43", "D44", "D46", "D47")
foreach ($cellRef in $textFormatCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.976.13'
$ws.Range("E2").Value = '  -0.24%  '
$ws.Range("D3").Value = '1.678.24'
$ws.Range("E3").Value = '  +0.18%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '215.10'
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("E6").Value = '  +1.34%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '0.251'
$ws.Range("E8").Value = '  -0.28%  '
$ws.Range("E9").Value = '  +0.38%  '
$ws.Range("D10").Value = '20.35'
$ws.Range("E10").Value = '  +0.69%  '
$ws.Range("E11").Value = '  -0.12%  '
$ws.Range("D12").Value = '1.914.23'
$ws.Range("E12").Value = '  +0.14%  '
$ws.Range("D13").Value = '1.648.48'
$ws.Range("E13").Value = '  -1.58%  '
$ws.Range("E14").Value = '  -0.02%  '
$ws.Range("D15").Value = '0.529'
$ws.Range("E15").Value = '  +1.12%  '
$ws.Range("D16").Value = '65.91'
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("D17").Value = '26.988.90'
$ws.Range("E17").Value = '  -0.30%  '
$ws.Range("D18").Value = '237.31'
$ws.Range("E18").Value = '  +0.26%  '
$ws.Range("D19").Value = '8.07'
$ws.Range("E19").Value = '  +4.05%  '
$ws.Range("E20").Value = '  -0.73%  '
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("E22").Value = '  -0.79%  '
$ws.Range("E23").Value = '  -0.95%  '
$ws.Range("E24").Value = '  -1.57%  '
$ws.Range("D25").Value = '145.59'
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("E26").Value = '  +1.46%  '
$ws.Range("D27").Value = '16.08'
$ws.Range("E27").Value = '  +0.68%  '
$ws.Range("E28").Value = '  -1.34%  '
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("D30").Value = '0.0498'
$ws.Range("E30").Value = '  -0.26%  '
$ws.Range("E31").Value = '  -0.50%  '
$ws.Range("E32").Value = '  +0.29%  '
$ws.Range("D33").Value = '1.488.07'
$ws.Range("E33").Value = '  +1.04%  '
$ws.Range("E34").Value = '  +0.67%  '
$ws.Range("E35").Value = '  +4.25%  '
$ws.Range("D36").Value = '2.42'
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("D37").Value = '0.586'
$ws.Range("E37").Value = '  +1.81%  '
$ws.Range("E38").Value = '  +2.75%  '
$ws.Range("D39").Value = '0.901'
$ws.Range("E39").Value = '  +0.04%  '
$ws.Range("E40").Value = '  -4.04%  '
$ws.Range("E41").Value = '  +0.42%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").Value = '2.32'
$ws.Range("E43").Value = '  +1.69%  '
$ws.Range("D44").Value = '67.59'
$ws.Range("E44").Value = '  +1.41%  '
$ws.Range("D45").Value = '1.819.92'
$ws.Range("E45").Value = '  -0.23%  '
$ws.Range("D46").Value = '0.779'
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("D47").Value = '90.57'
$ws.Range("E47").Value = '  +0.27%  '
$ws.Range("E48").Value = '  +15.55%  '
$ws.Range("E49").Value = '  -0.54%  '
$ws.Range("E50").Value = '  +1.68%  '
$ws.Range("E51").Value = '  +0.37%  '
